$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "43.223.85"
Set-TextValue "E2" "  -1.48%  "
Set-TextValue "D3" "2.268.59"
Set-TextValue "E3" "  -1.93%  "
Set-TextValue "E4" "  +0.23%  "
Set-TextValue "D5" "113.21"
Set-TextValue "E5" "  +4.18%  "
Set-TextValue "D6" "264.68"
Set-TextValue "E6" "  -2.64%  "
Set-TextValue "E7" "  -1.21%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.597"
Set-TextValue "E9" "  -3.54%  "
Set-TextValue "D10" "48.15"
Set-TextValue "E10" "  -0.08%  "
Set-TextValue "E11" "  -2.03%  "
Set-TextValue "D12" "8.72"
Set-TextValue "E12" "  +3.65%  "
Set-TextValue "E13" "  -0.44%  "
Set-TextValue "D14" "15.39"
Set-TextValue "E14" "  -2.74%  "
Set-TextValue "D15" "2.607.18"
Set-TextValue "E15" "  -1.91%  "
Set-TextValue "E16" "  -1.03%  "
Set-TextValue "D17" "2.267.89"
Set-TextValue "E17" "  -1.65%  "
Set-TextValue "D18" "43.117.98"
Set-TextValue "E18" "  -1.64%  "
Set-TextValue "E19" "  -3.87%  "
Set-TextValue "D20" "6.95"
Set-TextValue "E20" "  +9.95%  "
Set-TextValue "D21" "71.10"
Set-TextValue "E21" "  -1.83%  "
Set-TextValue "E22" "  -3.65%  "
Set-TextValue "D23" "9.81"
Set-TextValue "E23" "  +5.03%  "
Set-TextValue "D24" "230.18"
Set-TextValue "E24" "  -1.79%  "
Set-TextValue "E25" "  -4.46%  "
Set-TextValue "E26" "  -0.08%  "
Set-TextValue "E27" "  -0.90%  "
Set-TextValue "D28" "3.87"
Set-TextValue "E28" "  -2.01%  "
Set-TextValue "D29" "40.99"
Set-TextValue "E29" "  -0.14%  "
Set-TextValue "E30" "  -2.32%  "
Set-TextValue "E31" "  -1.69%  "
Set-TextValue "D32" "171.51"
Set-TextValue "E32" "  -3.35%  "
Set-TextValue "D33" "21.27"
Set-TextValue "E33" "  -3.13%  "
Set-TextValue "D34" "0.0904"
Set-TextValue "E34" "  -1.50%  "
Set-TextValue "D35" "5.62"
Set-TextValue "E35" "  +0.46%  "
Set-TextValue "E36" "  -0.57%  "
Set-TextValue "E37" "  -4.86%  "
Set-TextValue "D38" "0.0351"
Set-TextValue "E38" "  -2.16%  "
Set-TextValue "D39" "3.83"
Set-TextValue "E39" "  -1.31%  "
Set-TextValue "E40" "  -8.27%  "
Set-TextValue "D41" "14.20"
Set-TextValue "E41" "  +15.74%  "
Set-TextValue "D42" "74.82"
Set-TextValue "E42" "  +11.13%  "
Set-TextValue "D43" "2.44"
Set-TextValue "E43" "  +3.32%  "
Set-TextValue "D44" "0.235"
Set-TextValue "E44" "  -1.40%  "
Set-TextValue "D45" "6.10"
Set-TextValue "E45" "  +8.84%  "
Set-TextValue "E46" "  +0.09%  "
Set-TextValue "E47" "  -2.55%  "
Set-TextValue "E48" "  -2.46%  "
Set-TextValue "E49" "  -3.08%  "
Set-TextValue "D50" "100.48"
Set-TextValue "E50" "  +0.83%  "
Set-TextValue "E51" "  +0.41%  "
